# [pr4-2911] CutfillerCoefficient.Title remove from configuration
#
# The "Title" column (header "Title", single data value "Współczynnik") is
# dropped from the CutfillerCoefficient sheet/table - it's no longer part of
# the CutfillerCoefficient schema. CFTProductivityRateMin / CFTProductivityRateMax
# shift left into columns A/B. The CutfillerCoefficient sheet becomes the
# active tab (it was "Usage" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CutfillerCoefficient")

# Select & delete column A ("Title"/"Współczynnik"), shifting B,C left to A,B.
[void]$ws.Columns.Item(1).Select()
[void]$ws.Columns.Item(1).Delete()

# The sheet's XML-mapped table (Tabela6) still spans the old A1:C2 range with
# 3 columns; shrink it to the new A1:B2 / 2-column extent and make sure the
# column headers it tracks match what's actually in the header row now.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B2"))
$ws.Range("A1").Value = "CFTProductivityRateMin"
$ws.Range("B1").Value = "CFTProductivityRateMax"

# CutfillerCoefficient becomes the active sheet/tab.
$ws.Activate()
